$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edit the "TAMAÑO POR PORCIÓN" value for the Crema de leche row: remove the
# stray space in "15 mg" so it reads "15mg".
$ws.Range("I2").Value = "15mg"

# Reflect the resulting active selection as captured in the saved workbook.
$ws.Range("I2").Select()
